$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
$ws2 = $wb.Worksheets.Item("ARM")
$ws3 = $wb.Worksheets.Item("BSM")
$ws4 = $wb.Worksheets.Item("CRP")
$ws5 = $wb.Worksheets.Item("CUL")
$ws6 = $wb.Worksheets.Item("GSM")
$ws7 = $wb.Worksheets.Item("LTW")
$ws8 = $wb.Worksheets.Item("WVR")

# Sheet ALC row 8
$ws1.Range("H8").Value = 54.4
$ws1.Range("I8").Value = 55.25
$ws1.Range("J8").Value = 51
$ws1.Range("K8").Value = 165.75
$ws1.Range("L8").Value = 153
$ws1.Range("M8").Value = -26.75
$ws1.Range("N8").Value = -431

# Sheet ALC row 28
$ws1.Range("H28").Value = 1033.2858
$ws1.Range("I28").Value = 888.8333
$ws1.Range("K28").Value = 888.8333
$ws1.Range("M28").Value = -403.8333

# Sheet ALC row 98
$ws1.Range("H98").Value = 2499
$ws1.Range("I98").Value = 2499
$ws1.Range("K98").Value = 2499
$ws1.Range("M98").Value = -1001

# Sheet ALC row 122
$ws1.Range("H122").Value = 2499
$ws1.Range("I122").Value = 2499
$ws1.Range("K122").Value = 7497
$ws1.Range("M122").Value = -5047

# Sheet ALC row 135
$ws1.Range("H135").Value = 1593.1852
$ws1.Range("I135").Value = 558.9474
$ws1.Range("K135").Value = 5030.5266
$ws1.Range("M135").Value = -2495.5266

# Sheet ALC row 137
$ws1.Range("H137").Value = 957803.7
$ws1.Range("I137").Value = 2651.4
$ws1.Range("K137").Value = 7954.200000000001
$ws1.Range("M137").Value = -5404.200000000001

# Sheet ALC row 138
$ws1.Range("H138").Value = 1809.1708
$ws1.Range("I138").Value = 1074.3846
$ws1.Range("J138").Value = 3082.8
$ws1.Range("K138").Value = 3223.1538
$ws1.Range("L138").Value = 9248.400000000001
$ws1.Range("M138").Value = 1916.8462
$ws1.Range("N138").Value = -19528.4

# Sheet ARM row 2
$ws2.Range("H2").Value = 7505.391
$ws2.Range("I2").Value = 1764.4546
$ws2.Range("K2").Value = 1764.4546
$ws2.Range("M2").Value = -1651.4546

# Sheet ARM row 32
$ws2.Range("H32").Value = 7249944.5
$ws2.Range("I32").Value = 7755727
$ws2.Range("K32").Value = 7755727
$ws2.Range("M32").Value = -7755440

# Sheet ARM row 38
$ws2.Range("H38").Value = 0
$ws2.Range("I38").Value = 0
$ws2.Range("J38").Value = 0
$ws2.Range("K38").Value = 0
$ws2.Range("L38").Value = 0
$ws2.Range("M38").ClearContents()
$ws2.Range("N38").ClearContents()

# Sheet ARM row 116
$ws2.Range("H116").Value = 7505.391
$ws2.Range("I116").Value = 1764.4546
$ws2.Range("K116").Value = 1764.4546
$ws2.Range("M116").Value = 529.5454

# Sheet ARM row 132
$ws2.Range("H132").Value = 1484696.8
$ws2.Range("I132").Value = 1696224.9
$ws2.Range("K132").Value = 5088674.699999999
$ws2.Range("M132").Value = -5086144.699999999

# Sheet BSM row 3
$ws3.Range("H3").Value = 7505.391
$ws3.Range("I3").Value = 1764.4546
$ws3.Range("K3").Value = 1764.4546
$ws3.Range("M3").Value = -1650.4546

# Sheet BSM row 22
$ws3.Range("H22").Value = 477.9091
$ws3.Range("I22").Value = 437
$ws3.Range("K22").Value = 437
$ws3.Range("M22").Value = -264

# Sheet BSM row 51
$ws3.Range("H51").Value = 0
$ws3.Range("J51").Value = 0
$ws3.Range("L51").Value = 0
$ws3.Range("N51").ClearContents()

# Sheet BSM row 107
$ws3.Range("H107").Value = 4290.773
$ws3.Range("I107").Value = 3188.6155
$ws3.Range("K107").Value = 3188.6155
$ws3.Range("M107").Value = -1268.6155

# Sheet BSM row 134
$ws3.Range("H134").Value = 1094635.8
$ws3.Range("I134").Value = 1193456.8
$ws3.Range("J134").Value = 765232.3
$ws3.Range("K134").Value = 3580370.4
$ws3.Range("L134").Value = 2295696.9
$ws3.Range("M134").Value = -3577835.4
$ws3.Range("N134").Value = -2300766.9

# Sheet CRP row 19
$ws4.Range("H19").Value = 3024.7144
$ws4.Range("I19").Value = 5056
$ws4.Range("J19").Value = 316.33334
$ws4.Range("K19").Value = 5056
$ws4.Range("L19").Value = 316.33334
$ws4.Range("M19").Value = -4886
$ws4.Range("N19").Value = -656.33334

# Sheet CRP row 24
$ws4.Range("H24").Value = 3024.7144
$ws4.Range("I24").Value = 5056
$ws4.Range("J24").Value = 316.33334
$ws4.Range("K24").Value = 5056
$ws4.Range("L24").Value = 316.33334
$ws4.Range("M24").Value = -4886
$ws4.Range("N24").Value = -656.33334

# Sheet CRP row 31
$ws4.Range("H31").Value = 130512.96
$ws4.Range("I31").Value = 165516.8
$ws4.Range("K31").Value = 165516.8
$ws4.Range("M31").Value = -165221.8

# Sheet CRP row 34
$ws4.Range("H34").Value = 130512.96
$ws4.Range("I34").Value = 165516.8
$ws4.Range("K34").Value = 165516.8
$ws4.Range("M34").Value = -165314.8

# Sheet CRP row 35
$ws4.Range("H35").Value = 11793
$ws4.Range("I35").Value = 11793
$ws4.Range("K35").Value = 11793
$ws4.Range("M35").Value = -11499

# Sheet CRP row 86
$ws4.Range("H86").Value = 9453
$ws4.Range("I86").Value = 9980
$ws4.Range("K86").Value = 9980
$ws4.Range("M86").Value = -8857

# Sheet CRP row 89
$ws4.Range("H89").Value = 9453
$ws4.Range("I89").Value = 9980
$ws4.Range("K89").Value = 49900
$ws4.Range("M89").Value = -44284

# Sheet CRP row 105
$ws4.Range("H105").Value = 41405.89
$ws4.Range("I105").Value = 51993.57
$ws4.Range("K105").Value = 51993.57
$ws4.Range("M105").Value = -50246.57

# Sheet CUL row 23
$ws5.Range("H23").Value = 107.875
$ws5.Range("I23").Value = 48.833332
$ws5.Range("K23").Value = 146.499996
$ws5.Range("M23").Value = 88.50000399999999

# Sheet CUL row 32
$ws5.Range("H32").Value = 494845250
$ws5.Range("J32").Value = 976190460
$ws5.Range("L32").Value = 2928571380
$ws5.Range("N32").Value = -2928571946

# Sheet CUL row 54
$ws5.Range("H54").Value = 0
$ws5.Range("I54").Value = 0
$ws5.Range("K54").Value = 0
$ws5.Range("M54").ClearContents()

# Sheet CUL row 63
$ws5.Range("H63").Value = 12611.2
$ws5.Range("I63").Value = 3000
$ws5.Range("K63").Value = 9000
$ws5.Range("M63").Value = -8251

# Sheet CUL row 66
$ws5.Range("H66").Value = 12611.2
$ws5.Range("I66").Value = 3000
$ws5.Range("K66").Value = 27000
$ws5.Range("M66").Value = -23256

# Sheet CUL row 75
$ws5.Range("H75").Value = 5163.5
$ws5.Range("J75").Value = 5961.154
$ws5.Range("L75").Value = 17883.462
$ws5.Range("N75").Value = -19879.462

# Sheet CUL row 78
$ws5.Range("H78").Value = 5163.5
$ws5.Range("J78").Value = 5961.154
$ws5.Range("L78").Value = 53650.38600000001
$ws5.Range("N78").Value = -63634.38600000001

# Sheet CUL row 81
$ws5.Range("H81").Value = 7658.6
$ws5.Range("J81").Value = 6323.25
$ws5.Range("L81").Value = 18969.75
$ws5.Range("N81").Value = -21215.75

# Sheet CUL row 84
$ws5.Range("H84").Value = 7658.6
$ws5.Range("J84").Value = 6323.25
$ws5.Range("L84").Value = 56909.25
$ws5.Range("N84").Value = -68141.25

# Sheet CUL row 108
$ws5.Range("H108").Value = 12127.444
$ws5.Range("I108").Value = 3342.3333
$ws5.Range("K108").Value = 10026.9999
$ws5.Range("M108").Value = -7146.999899999999

# Sheet CUL row 132
$ws5.Range("H132").Value = 2984.8333
$ws5.Range("I132").Value = 1402.2858
$ws5.Range("J132").Value = 5200.4
$ws5.Range("K132").Value = 12620.5722
$ws5.Range("L132").Value = 46803.6
$ws5.Range("M132").Value = -10090.5722
$ws5.Range("N132").Value = -51863.6

# Sheet CUL row 140
$ws5.Range("H140").Value = 2663.394
$ws5.Range("I140").Value = 1686
$ws5.Range("K140").Value = 5058
$ws5.Range("M140").Value = 122

# Sheet GSM row 18
$ws6.Range("H18").Value = 4989.25
$ws6.Range("I18").Value = 3319
$ws6.Range("J18").Value = 10000
$ws6.Range("K18").Value = 3319
$ws6.Range("L18").Value = 10000
$ws6.Range("M18").Value = -3026
$ws6.Range("N18").Value = -10586

# Sheet GSM row 70
$ws6.Range("H70").Value = 4656.5
$ws6.Range("I70").Value = 4662.2666
$ws6.Range("J70").Value = 4627.6665
$ws6.Range("K70").Value = 4662.2666
$ws6.Range("L70").Value = 4627.6665
$ws6.Range("M70").Value = -4392.2666
$ws6.Range("N70").Value = -5167.6665

# Sheet GSM row 73
$ws6.Range("H73").Value = 4656.5
$ws6.Range("I73").Value = 4662.2666
$ws6.Range("J73").Value = 4627.6665
$ws6.Range("K73").Value = 4662.2666
$ws6.Range("L73").Value = 4627.6665
$ws6.Range("M73").Value = -3726.2666
$ws6.Range("N73").Value = -6499.6665

# Sheet GSM row 110
$ws6.Range("H110").Value = 99508.8
$ws6.Range("J110").Value = 99508.8
$ws6.Range("L110").Value = 99508.8
$ws6.Range("N110").Value = -107688.8

# Sheet GSM row 113
$ws6.Range("H113").Value = 0
$ws6.Range("I113").Value = 0
$ws6.Range("J113").Value = 0
$ws6.Range("K113").Value = 0
$ws6.Range("L113").Value = 0
$ws6.Range("M113").ClearContents()
$ws6.Range("N113").ClearContents()

# Sheet LTW row 122
$ws7.Range("H122").Value = 5399.0557
$ws7.Range("I122").Value = 5175.3335
$ws7.Range("J122").Value = 5846.5
$ws7.Range("K122").Value = 15526.0005
$ws7.Range("L122").Value = 17539.5
$ws7.Range("M122").Value = -13076.0005
$ws7.Range("N122").Value = -22439.5

# Sheet LTW row 132
$ws7.Range("H132").Value = 4354836
$ws7.Range("I132").Value = 11600896
$ws7.Range("J132").Value = 7199.6
$ws7.Range("K132").Value = 34802688
$ws7.Range("L132").Value = 21598.8
$ws7.Range("M132").Value = -34800158
$ws7.Range("N132").Value = -26658.8

# Sheet WVR row 82
$ws8.Range("H82").Value = 0
$ws8.Range("J82").Value = 0
$ws8.Range("L82").Value = 0
$ws8.Range("N82").ClearContents()

# Sheet WVR row 85
$ws8.Range("H85").Value = 0
$ws8.Range("J85").Value = 0
$ws8.Range("L85").Value = 0
$ws8.Range("N85").ClearContents()

# Sheet WVR row 96
$ws8.Range("H96").Value = 2792.2856
$ws8.Range("I96").Value = 2423
$ws8.Range("J96").Value = 3284.6667
$ws8.Range("K96").Value = 2423
$ws8.Range("L96").Value = 3284.6667
$ws8.Range("M96").Value = -1050
$ws8.Range("N96").Value = -6030.6667

# Sheet WVR row 132
$ws8.Range("H132").Value = 5922285
$ws8.Range("I132").Value = 6101293.5
$ws8.Range("K132").Value = 18303880.5
$ws8.Range("M132").Value = -18301350.5

# Sheet WVR row 136
$ws8.Range("H136").Value = 10244.826
$ws8.Range("J136").Value = 7525.3335
$ws8.Range("L136").Value = 22576.0005
$ws8.Range("N136").Value = -27676.0005
